$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing "extr" rows (rows 8-15): columns C, D, E ---
$ws.Range("C8").Value = 14
$ws.Range("D8").Value = 11

$ws.Range("C9").Value = 16
$ws.Range("E9").Value = $true

$ws.Range("C10").Value = 5
$ws.Range("D10").Value = 12
$ws.Range("E10").Value = $true

$ws.Range("C11").Value = 5
$ws.Range("D11").Value = 9

$ws.Range("C12").Value = 10
$ws.Range("E12").Value = $true

$ws.Range("D13").Value = 8
$ws.Range("E13").Value = $true

$ws.Range("C14").Value = 9
$ws.Range("D14").Value = 11
$ws.Range("E14").Value = $true

$ws.Range("C15").Value = 7
$ws.Range("D15").Value = 11

# --- Append two new rows (16, 17) for line7 / line8 ---
# Copy formatting from the row above down into the new rows first, then set values,
# so the new cells pick up the same style (bordered/bold/centered) as A2:A15.
$ws.Range("A15").Copy($ws.Range("A16"))
$ws.Range("A15").Copy($ws.Range("A17"))

$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "line7"
$ws.Range("C16").Value = 5
$ws.Range("D16").Value = 7
$ws.Range("E16").Value = $false

$ws.Range("A17").Value = 15
$ws.Range("B17").Value = "line8"
$ws.Range("C17").Value = 8
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = $true
